# This script updates the "Stationary generator alpha zero" experiment
# workbook with a new (x, y) evaluation point and all of the dependent
# values that flow from it (restrictions, vector_bf, vector_BF, etc.).
#
# Many of the target cells originally hold numeric-looking values that
# were stored as *text* (shared-string) cells, not as real numbers.
# Assigning a numeric-looking string straight to Range.Value causes
# Excel to auto-convert it to a genuine number cell, so for those cells
# we use the classic "leading apostrophe" trick to force text entry and
# then restore the Normal style so no stray number formatting sticks
# around on the cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Forces a cell to contain a literal text value even when the text
    # looks like a number (e.g. "-3.3000000000000003").
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# NOTE: worksheet lookup by name (Worksheets.Item("Name")) is
# case-insensitive in this environment, and this workbook has two
# sheets whose names differ only by case ("Vector_bf" / "Vector_BF").
# To avoid ambiguity we always address sheets by their (1-based)
# position, which matches the <sheets> order in workbook.xml:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---------------------------------------------------------------
# Restricciones_del_lider
# ---------------------------------------------------------------
$wsLider = $wb.Worksheets.Item(2)

$wsLider.Range("A2").Value = "2.3000000000000003 - x"
Set-TextValue $wsLider.Range("B2") "-3.3000000000000003"
Set-TextValue $wsLider.Range("D2") "0.51"

$wsLider.Range("A3").Value = "-2.3000000000000003 + x"
Set-TextValue $wsLider.Range("B3") "1.3000000000000003"
Set-TextValue $wsLider.Range("D3") "0.17"

# ---------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item(3)

$wsFollower.Range("A2").Value = "-4.449999999999999 + y"
Set-TextValue $wsFollower.Range("B2") "3.4499999999999993"
Set-TextValue $wsFollower.Range("D2") "0.82"
Set-TextValue $wsFollower.Range("E2") "1.7000000000000002"
Set-TextValue $wsFollower.Range("F2") "8.5"

$wsFollower.Range("A3").Value = "4.449999999999999 - y"
Set-TextValue $wsFollower.Range("B3") "-5.449999999999999"
Set-TextValue $wsFollower.Range("D3") "0.81"
Set-TextValue $wsFollower.Range("E3") "0"
Set-TextValue $wsFollower.Range("F3") "5.1"

# ---------------------------------------------------------------
# Punto_modificado (x, y)
# ---------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)

Set-TextValue $wsPunto.Range("A2") "2.3000000000000003"
Set-TextValue $wsPunto.Range("B2") "4.449999999999999"

# ---------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------
$wsVecBf = $wb.Worksheets.Item(5)

Set-TextValue $wsVecBf.Range("A2") "1.9220000000000041"

# ---------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------
$wsVecBF = $wb.Worksheets.Item(6)

Set-TextValue $wsVecBF.Range("A2") "-0.66"
Set-TextValue $wsVecBF.Range("A3") "-0.7000000000000002"
